$wb = $excel.ActiveWorkbook

# Sheet "展览": update F2 (645 -> 646) and F4 (1462 -> 1468)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 646
$ws1.Range("F4").Value = 1468

# Sheet "全部类型": update F2 (645 -> 646) and F4 (1462 -> 1468)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 646
$ws4.Range("F4").Value = 1468
